# Realizar Receita (Ingredientes) - align with "Selecionar Receita" use case
# as described by Petra: drop the separate "«include» Selecionar Receita"
# step (old row 7) and renumber the remaining steps (shift 3.x -> 2.x, 4 -> 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old row 7 ("1. «include» Selecionar Receita") - this shifts
#    every following row up by one and keeps all number formats / merges in
#    sync automatically (mirrors what Excel does for a row delete).
$ws.Rows.Item(7).Delete()

# 2) Row height tweaks that came with the edit (row 7 is now the
#    "description" row, row 8 became the "quer iniciar" row).
$ws.Rows.Item(7).RowHeight = 35.25
$ws.Rows.Item(8).RowHeight = 19.5

# 3) Pre-condition text changed.
$ws.Range("C4").Value = "Ter Selecionado Receita"

# 4) Renumber the remaining scenario steps (3.x -> 2.x, 4. -> 3., etc).
$ws.Range("D7").Value = "1. Mostra Descrição da Receita e respetivos Ingredientes"
$ws.Range("C8").Value = "2. Informa que quer Iniciar a Receita"
$ws.Range("D9").Value = "3. Inicia a Receita"

$ws.Range("C10").Value = "2.1. Informa que faltam Ingredientes"
$ws.Range("B11").Value = "[Falta ingredientes] (passo 2)"
$ws.Range("D11").Value = "2.2. Obtém localização Utilizador e calcula supermercados mais próximos"
$ws.Range("D12").Value = "2.3. Mostra supermecados mais próximos"
$ws.Range("C13").Value = "2.4. Informa que pode prosseguir com a Receita"
$ws.Range("D14").Value = "2.5. Volta a 4."

$ws.Range("C15").Value = "2.1. Informa que não quer realizar Receita"
$ws.Range("D16").Value = "2.2. Cancela Realização de Receita"
$ws.Range("B17").Value = "(passo 2)"

$ws.Range("C18").Value = "2.4.1. Informa que quer cancelar Receita"
$ws.Range("D19").Value = "2.4.2. Cancela Realização de Receita"
$ws.Range("B20").Value = "(passo 2.4)"

# 5) Selection moved to B17 in the saved file.
$ws.Range("B17").Select()
